# 496680 Added comment to US114 - moved to Build 2
#
# US114 is row 4 of Sheet1 (A4 = "US114"). Add a note in the blank
# "Rational" column cell (F4) saying the story moved to Build 2, and
# leave the selection sitting on that cell, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cell = $ws.Range("F4")
$cell.Value = "MOVED TO BUILD 2"

# Left-align horizontally, keep/center vertically (matches the rest of
# the data rows), same default font/fill/border as an un-styled cell.
$cell.HorizontalAlignment = -4131   # xlLeft
$cell.VerticalAlignment = -4108     # xlCenter

$ws.Range("F4").Select()
